$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (test_number) values from 2 to 3 for rows 2-11
$ws.Range("A2:A11").Value = 3

# Update the active cell / selection to I9
$ws.Range("I9").Select()
